# Heap Implementation containing all the methods of HEAP ADT like
# percolateUp, percolateDown — add a new "MedianInfiniteStream" class-name
# entry to the table for the "Median of numbers in infinite stream of
# integers" row (row 5), styled like the other ClassName column entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New cell H5: ClassName for the row-5 problem.
$cell = $ws.Range("H5")
$cell.Value = "MedianInfiniteStream"
$cell.Font.Name = "Comic Sans MS"
$cell.Font.Size = 9.8
$cell.VerticalAlignment = -4108   # xlCenter (vertical centering)

# Row 5 grows a bit taller to accommodate the new font.
$ws.Rows.Item(5).RowHeight = 15.6

# Move/save the active selection on the sheet to the newly-filled cell.
$cell.Select()

# Set the sheet to print in portrait orientation.
$ws.PageSetup.Orientation = 1   # xlPortrait
